$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Selection Sort" / "Insertion Sort" column headers (A1 <-> B1).
$ws.Range("A1").Value = "Insertion Sort"
$ws.Range("B1").Value = "Selection Sort"

# Reposition / resize the embedded chart (it was dragged to a new size/position).
$co = $ws.ChartObjects().Item(1)
$co.Top = 0.6
$co.Left = 403.3294921875
$co.Width = 642.5802734375
$co.Height = 409.5

# Restore the last active selection.
$ws.Range("O5").Select()
